$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 <- old Row 21
$ws.Range("A3").Value2 = 111866170
$ws.Range("B3").Value2 = 90682
$ws.Range("E3").Value2 = 2059
$ws.Range('F3').Value2 = 'Skrovlig taggsvamp'
$ws.Range('G3').Value2 = 'Hydnellum scabrosum'
$ws.Range('H3').Value2 = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Range("Q3").Value2 = 702754.3208386695
$ws.Range("R3").Value2 = 7299886.818591502

# Row 4 <- old Row 19
$ws.Range("A4").Value2 = 111865981
$ws.Range("B4").Value2 = 90652
$ws.Range("E4").Value2 = 3100
$ws.Range('F4').Value2 = 'Talltaggsvamp'
$ws.Range('G4').Value2 = 'Bankera fuligineoalba'
$ws.Range('H4').Value2 = '(Schmidt : Fr.) Pouzar'
$ws.Range("Q4").Value2 = 702695.6801449896
$ws.Range("R4").Value2 = 7299770.100652335

# Row 5 <- old Row 9
$ws.Range("A5").Value2 = 111866131
$ws.Range("B5").Value2 = 90682
$ws.Range("E5").Value2 = 2059
$ws.Range('F5').Value2 = 'Skrovlig taggsvamp'
$ws.Range('G5').Value2 = 'Hydnellum scabrosum'
$ws.Range('H5').Value2 = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Range("Q5").Value2 = 702756.5806601554
$ws.Range("R5").Value2 = 7299854.813386399
$ws.Range('AC5').Value2 = 'Flera fruktkoppar som växer i en häxring'

# Row 6 <- old Row 5
$ws.Range("A6").Value2 = 111866252
$ws.Range("B6").Value2 = 78107
$ws.Range("E6").Value2 = 6453
$ws.Range('F6').Value2 = 'Vedskivlav'
$ws.Range('G6').Value2 = 'Hertelidea botryosa'
$ws.Range('H6').Value2 = '(Fr.) Printzen & Kantvilas'
$ws.Range("Q6").Value2 = 702680.6244306123
$ws.Range("R6").Value2 = 7299924.914052285

# Row 7 <- old Row 18
$ws.Range("A7").Value2 = 111866021
$ws.Range("B7").Value2 = 78107
$ws.Range("E7").Value2 = 6453
$ws.Range('F7').Value2 = 'Vedskivlav'
$ws.Range('G7').Value2 = 'Hertelidea botryosa'
$ws.Range('H7').Value2 = '(Fr.) Printzen & Kantvilas'
$ws.Range("Q7").Value2 = 702738.1111920479
$ws.Range("R7").Value2 = 7299806.49869829

# Row 8 <- old Row 3
$ws.Range("A8").Value2 = 111866159
$ws.Range("B8").Value2 = 90652
$ws.Range("E8").Value2 = 3100
$ws.Range('F8').Value2 = 'Talltaggsvamp'
$ws.Range('G8').Value2 = 'Bankera fuligineoalba'
$ws.Range('H8').Value2 = '(Schmidt : Fr.) Pouzar'
$ws.Range("Q8").Value2 = 702755.4455659754
$ws.Range("R8").Value2 = 7299865.042498757

# Row 9 <- old Row 6
$ws.Range("A9").Value2 = 111866301
$ws.Range("B9").Value2 = 90660
$ws.Range("E9").Value2 = 4362
$ws.Range('F9').Value2 = 'Blå taggsvamp'
$ws.Range('G9').Value2 = 'Hydnellum caeruleum'
$ws.Range('H9').Value2 = '(Hornem.) P.Karst.'
$ws.Range("Q9").Value2 = 702522.1051459431
$ws.Range("R9").Value2 = 7300047.742725079
$ws.Range("AC9").ClearContents()

# Row 10 <- old Row 4
$ws.Range("A10").Value2 = 111866276
$ws.Range("B10").Value2 = 78107
$ws.Range("E10").Value2 = 6453
$ws.Range('F10').Value2 = 'Vedskivlav'
$ws.Range('G10').Value2 = 'Hertelidea botryosa'
$ws.Range('H10').Value2 = '(Fr.) Printzen & Kantvilas'
$ws.Range("Q10").Value2 = 702660.5304515015
$ws.Range("R10").Value2 = 7299928.856484808

# Row 11 <- old Row 22
$ws.Range("A11").Value2 = 111866194
$ws.Range("B11").Value2 = 90682
$ws.Range("E11").Value2 = 2059
$ws.Range('F11').Value2 = 'Skrovlig taggsvamp'
$ws.Range('G11').Value2 = 'Hydnellum scabrosum'
$ws.Range('H11').Value2 = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Range("Q11").Value2 = 702686.7518818546
$ws.Range("R11").Value2 = 7299919.985876646
$ws.Range('AC11').Value2 = 'Flera fruktkoppar som växer i en häxring'

# Row 12 <- old Row 13
$ws.Range("A12").Value2 = 111866065
$ws.Range("B12").Value2 = 78107
$ws.Range("E12").Value2 = 6453
$ws.Range('F12').Value2 = 'Vedskivlav'
$ws.Range('G12').Value2 = 'Hertelidea botryosa'
$ws.Range('H12').Value2 = '(Fr.) Printzen & Kantvilas'
$ws.Range("Q12").Value2 = 702767.9701038125
$ws.Range("R12").Value2 = 7299827.988589783

# Row 13 <- old Row 11
$ws.Range("A13").Value2 = 111865488
$ws.Range("B13").Value2 = 90660
$ws.Range("E13").Value2 = 4362
$ws.Range('F13').Value2 = 'Blå taggsvamp'
$ws.Range('G13').Value2 = 'Hydnellum caeruleum'
$ws.Range('H13').Value2 = '(Hornem.) P.Karst.'
$ws.Range("Q13").Value2 = 702716.2360189059
$ws.Range("R13").Value2 = 7299724.539719297

# Row 14 <- old Row 10
$ws.Range("A14").Value2 = 111866048
$ws.Range("B14").Value2 = 90682
$ws.Range("E14").Value2 = 2059
$ws.Range('F14').Value2 = 'Skrovlig taggsvamp'
$ws.Range('G14').Value2 = 'Hydnellum scabrosum'
$ws.Range('H14').Value2 = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'

# Row 15 <- old Row 14
$ws.Range("A15").Value2 = 111866031
$ws.Range("B15").Value2 = 78107
$ws.Range('D15').Value2 = 'NT'
$ws.Range("E15").Value2 = 6453
$ws.Range('F15').Value2 = 'Vedskivlav'
$ws.Range('G15').Value2 = 'Hertelidea botryosa'
$ws.Range('H15').Value2 = '(Fr.) Printzen & Kantvilas'
$ws.Range("Q15").Value2 = 702750.1350314748
$ws.Range("R15").Value2 = 7299799.924799141
$ws.Range("L15").ClearContents()

# Row 16 <- old Row 7
$ws.Range("A16").Value2 = 111865866
$ws.Range("B16").Value2 = 90652
$ws.Range("E16").Value2 = 3100
$ws.Range('F16').Value2 = 'Talltaggsvamp'
$ws.Range('G16').Value2 = 'Bankera fuligineoalba'
$ws.Range('H16').Value2 = '(Schmidt : Fr.) Pouzar'
$ws.Range("Q16").Value2 = 702753.3055412351
$ws.Range("R16").Value2 = 7299801.798166115

# Row 17 <- old Row 8
$ws.Range("A17").Value2 = 111865263
$ws.Range("B17").Value2 = 90658
$ws.Range("E17").Value2 = 4361
$ws.Range('F17').Value2 = 'Orange taggsvamp'
$ws.Range('G17').Value2 = 'Hydnellum aurantiacum'
$ws.Range('H17').Value2 = '(Batsch:Fr.) P.Karst.'
$ws.Range("Q17").Value2 = 702714.1819675351
$ws.Range("R17").Value2 = 7299724.394724619

# Row 18 <- old Row 15
$ws.Range("A18").Value2 = 111865919
$ws.Range("B18").Value2 = 95538
$ws.Range('D18').Value2 = 'LC'
$ws.Range("E18").Value2 = 221941
$ws.Range('F18').Value2 = 'Plattlummer'
$ws.Range('G18').Value2 = 'Lycopodium complanatum'
$ws.Range('H18').Value2 = 'L.'
$ws.Range("Q18").Value2 = 702755.0230470664
$ws.Range("R18").Value2 = 7299754.083126943
$ws.Range("L18").Value2 = ""  # EMPTY_STR_MARKER

# Row 19 <- old Row 17
$ws.Range("A19").Value2 = 111865961
$ws.Range("B19").Value2 = 77267
$ws.Range("E19").Value2 = 6446
$ws.Range('F19').Value2 = 'Kolflarnlav'
$ws.Range('G19').Value2 = 'Carbonicola anthracophila'
$ws.Range('H19').Value2 = '(Nyl.) Bendiksby & Timdal'
$ws.Range("Q19").Value2 = 702714.4770808229
$ws.Range("R19").Value2 = 7299790.39698876

# Row 20 <- old Row 12
$ws.Range("A20").Value2 = 111865578
$ws.Range("B20").Value2 = 90854
$ws.Range("E20").Value2 = 2079
$ws.Range('F20').Value2 = 'Nordtagging'
$ws.Range('G20').Value2 = 'Odonticium romellii'
$ws.Range('H20').Value2 = '(S.Lundell) Parmasto'
$ws.Range("Q20").Value2 = 702741.9879008483
$ws.Range("R20").Value2 = 7299745.739876431

# Row 21 <- old Row 20
$ws.Range("A21").Value2 = 111865524
$ws.Range("B21").Value2 = 90660
$ws.Range("E21").Value2 = 4362
$ws.Range('F21').Value2 = 'Blå taggsvamp'
$ws.Range('G21').Value2 = 'Hydnellum caeruleum'
$ws.Range('H21').Value2 = '(Hornem.) P.Karst.'
$ws.Range("Q21").Value2 = 702731.0699128226
$ws.Range("R21").Value2 = 7299742.494774668

# Row 22 <- old Row 16
$ws.Range("A22").Value2 = 111865668
$ws.Range("B22").Value2 = 78107
$ws.Range("E22").Value2 = 6453
$ws.Range('F22').Value2 = 'Vedskivlav'
$ws.Range('G22').Value2 = 'Hertelidea botryosa'
$ws.Range('H22').Value2 = '(Fr.) Printzen & Kantvilas'
$ws.Range("Q22").Value2 = 702740.9003275807
$ws.Range("R22").Value2 = 7299743.601162716
$ws.Range("AC22").ClearContents()
